$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 382, shifting the existing rows 382-434 down to 386-438.
$ws.Rows("382:385").Insert()

# Common/boilerplate values shared by all rows in this data block.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

# New row 382: Calidad "Primera", Variedad stays "Cultivar IV Región"
$ws.Cells.Item(382, 1).Value = $mercadoId
$ws.Cells.Item(382, 2).Value = $mercado
$ws.Cells.Item(382, 3).Value = $region
$ws.Cells.Item(382, 4).Value = 45124
$ws.Cells.Item(382, 5).Value = $codreg
$ws.Cells.Item(382, 6).Value = $categoriaId
$ws.Cells.Item(382, 7).Value = $categoria
$ws.Cells.Item(382, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 420
$ws.Cells.Item(382, 11).Value = 15000
$ws.Cells.Item(382, 12).Value = 15000
$ws.Cells.Item(382, 13).Value = 15000
$ws.Cells.Item(382, 14).Value = $unidad
$ws.Cells.Item(382, 15).Value = $origen
$ws.Cells.Item(382, 16).Value = 833
$ws.Cells.Item(382, 17).Value = $kgUnidades
$ws.Cells.Item(382, 18).Value = $clasificacion

# New row 383: Calidad "Extra", Variedad "Sin especificar"
$ws.Cells.Item(383, 1).Value = $mercadoId
$ws.Cells.Item(383, 2).Value = $mercado
$ws.Cells.Item(383, 3).Value = $region
$ws.Cells.Item(383, 4).Value = 45124
$ws.Cells.Item(383, 5).Value = $codreg
$ws.Cells.Item(383, 6).Value = $categoriaId
$ws.Cells.Item(383, 7).Value = $categoria
$ws.Cells.Item(383, 8).Value = "Sin especificar"
$ws.Cells.Item(383, 9).Value = "Extra"
$ws.Cells.Item(383, 10).Value = 270
$ws.Cells.Item(383, 11).Value = 16000
$ws.Cells.Item(383, 12).Value = 16000
$ws.Cells.Item(383, 13).Value = 16000
$ws.Cells.Item(383, 14).Value = $unidad
$ws.Cells.Item(383, 15).Value = $origen
$ws.Cells.Item(383, 16).Value = 889
$ws.Cells.Item(383, 17).Value = $kgUnidades
$ws.Cells.Item(383, 18).Value = $clasificacion

# New row 384: Calidad "Segunda", Variedad "Sin especificar"
$ws.Cells.Item(384, 1).Value = $mercadoId
$ws.Cells.Item(384, 2).Value = $mercado
$ws.Cells.Item(384, 3).Value = $region
$ws.Cells.Item(384, 4).Value = 45124
$ws.Cells.Item(384, 5).Value = $codreg
$ws.Cells.Item(384, 6).Value = $categoriaId
$ws.Cells.Item(384, 7).Value = $categoria
$ws.Cells.Item(384, 8).Value = "Sin especificar"
$ws.Cells.Item(384, 9).Value = "Segunda"
$ws.Cells.Item(384, 10).Value = 290
$ws.Cells.Item(384, 11).Value = 12000
$ws.Cells.Item(384, 12).Value = 12000
$ws.Cells.Item(384, 13).Value = 12000
$ws.Cells.Item(384, 14).Value = $unidad
$ws.Cells.Item(384, 15).Value = $origen
$ws.Cells.Item(384, 16).Value = 667
$ws.Cells.Item(384, 17).Value = $kgUnidades
$ws.Cells.Item(384, 18).Value = $clasificacion

# New row 385: Calidad "Tercera", Variedad "Sin especificar"
$ws.Cells.Item(385, 1).Value = $mercadoId
$ws.Cells.Item(385, 2).Value = $mercado
$ws.Cells.Item(385, 3).Value = $region
$ws.Cells.Item(385, 4).Value = 45124
$ws.Cells.Item(385, 5).Value = $codreg
$ws.Cells.Item(385, 6).Value = $categoriaId
$ws.Cells.Item(385, 7).Value = $categoria
$ws.Cells.Item(385, 8).Value = "Sin especificar"
$ws.Cells.Item(385, 9).Value = "Tercera"
$ws.Cells.Item(385, 10).Value = 230
$ws.Cells.Item(385, 11).Value = 10000
$ws.Cells.Item(385, 12).Value = 10000
$ws.Cells.Item(385, 13).Value = 10000
$ws.Cells.Item(385, 14).Value = $unidad
$ws.Cells.Item(385, 15).Value = $origen
$ws.Cells.Item(385, 16).Value = 556
$ws.Cells.Item(385, 17).Value = $kgUnidades
$ws.Cells.Item(385, 18).Value = $clasificacion
